$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = New-Object "object[,]" 154,2
$data[0,0] = "Clinical data not provided"
$data[0,1] = "TBT"
$data[1,0] = "Clinical Outcome"
$data[1,1] = "លទ្ធផលគ្លីនិក"
$data[2,0] = "Clinical Outcome Status:"
$data[2,1] = "ស្ថានភាពលទ្ធផលគ្លីនិក៖"
$data[3,0] = "Co-resistances"
$data[3,1] = "សហ-ភាពសុំា"
$data[4,0] = "Combine Susceptible + Intermediate"
$data[4,1] = "រួមបញ្ជូលគ្នានៃ Susceptible + Intermediate"
$data[5,0] = "Consider saving .acorn file on the cloud for additional security."
$data[5,1] = "ពិចារណាក្នុងការរក្សាទុកឯកសារ .acorn ក្នុង cloud សម្រាប់ការបន្ថែមសុវត្ថិភាព"
$data[6,0] = "Contains names of organisms before and after mapping."
$data[6,1] = "TBT"
$data[7,0] = "Couldn't connect to server. Please check internet access."
$data[7,1] = "មិនអាចភ្ជាប់ម៉ាស៊ីនមេទៅបានទេ។ សូមពិនិត្យមើលដំណើរការអ៊ីនធឺណិត។"
$data[8,0] = "Critical errors with clinical data."
$data[8,1] = "កំហុសឆ្គងសំខាន់ៗក្នុងផ្នែកទិន្នន័យគ្លីនិក។"
$data[9,0] = "Culture results per specimen type"
$data[9,1] = "លទ្ធផលការបណ្តុះមេរោគតាមប្រភេទវត្ថុវិភាគ"
$data[10,0] = "Data Management"
$data[10,1] = "ការគ្រប់គ្រងទិន្នន័យ"
$data[11,0] = "Date of Enrolment"
$data[11,1] = "ថ្ងៃចូលរួមការអង្កេតតាមដាន"
$data[12,0] = "Day 28"
$data[12,1] = "ថ្ងៃទី២៨"
$data[13,0] = "Day 28 Status:"
$data[13,1] = "ស្ថានភាពថ្ងៃទី២៨៖"
$data[14,0] = "Diagnosis at Enrolment"
$data[14,1] = "រោគវិនិច្ឆ័យនៅថ្ងៃចូលរួមការអង្កេតតាមដាន"
$data[15,0] = "Dismiss"
$data[15,1] = "ច្រានចោល"
$data[16,0] = "Distribution of Enrolments"
$data[16,1] = "របាយអ្នកចូលរួមការអង្កេតតាមដាន"
$data[17,0] = "Download Enrolment Log (.xlsx)"
$data[17,1] = "ទាញយកបញ្ជីអ្នកចូលរួមការអង្កេតតាមដាន (.xlsx)"
$data[18,0] = "Download Lab Log (.xlsx)"
$data[18,1] = "TBT"
$data[19,0] = "Empiric Antibiotics Prescribed"
$data[19,1] = "អង់ទីប៊ីយ៉ូទិកដែលផ្តល់"
$data[20,0] = "Enrolments"
$data[20,1] = "ការចូលរួមអង្កេតតាមដាន"
$data[21,0] = "Enrolments by (type of) Ward"
$data[21,1] = "ការចូលរួមការអង្កេតតាមដានតាម​ (ប្រភេទ) អគារ"
$data[22,0] = "Enrolments with Blood Culture"
$data[22,1] = "ការចូលរួមការអង្កេតតាមដានជាមួយការបណ្តុះមេរោគក្នុងឈាម"
$data[23,0] = "Error in combining clinical and lab data."
$data[23,1] = "មានកំហុសក្នុងការរួមបញ្ចូលទិន្នន័យគ្លីនិកនិងមន្ទីរពិសោធន៍។"
$data[24,0] = "Every D28 record (F04) matches exactly one patient enrolment (F01)."
$data[24,1] = "កំណត់ត្រា D28 (F04) ទាំងអស់ ត្រូវគ្នាជាមួយនឹងអ្នកជំងឺដែលចូលរួមការអង្កេតតាមដាន (F01)។"
$data[25,0] = "Every hospital outcome record (F03) has a matching infection episode (F02)."
$data[25,1] = "កំណត់ត្រាលទ្ធផលនៅមន្ទីរពេទ្យទាំងអស់(F03) ត្រូវគ្នាជាមួយ វគ្គការបង្ករោគ (F02)។"
$data[26,0] = "Every hospital outcome record (F03) has a matching patient enrolment (F01)."
$data[26,1] = "កំណត់ត្រាលទ្ធផលនៅមន្ទីរពេទ្យទាំងអស់ (F03) ត្រូវគ្នាជាមួយនឹងអ្នកជំងឺដែលចូលរួមការអង្កេតតាមដាន(F01)។"
$data[27,0] = "Every infection episode record (F02) has a matching patient enrolment (F01)."
$data[27,1] = "កំណត់ត្រាវគ្គការបង្ករោគទាំងអស់ (F02) ត្រូវគ្នាជាមួយនឹងអ្នកជំងឺដែលចូលរួមការអង្កេតតាមដាន (F01)។"
$data[28,0] = "File name:"
$data[28,1] = "ឈ្មោះឯកសារ៖"
$data[29,0] = "First sheet is the log of all enrolments retrived from REDCap (as per adjacent table). The second sheet is a listing of all flagged elements."
$data[29,1] = "សន្លឹកទីមួយគឺជាបញ្ជីអ្នកចូលរួមការអង្កេតតាមដានទាំងអស់ដែលបានទាញយកពី REDCap (ដូចទៅនឹងតារាងក្បែរគ្នា)។ សន្លឹកទីពីរគឺជាបញ្ជីនៃកំហុសឆ្គងទាំងអស់។"
$data[30,0] = "Follow-up"
$data[30,1] = "ការតាមដានបន្ត"
$data[31,0] = "from cultures that have growth"
$data[31,1] = "ពីការបណ្តុះមេរោគដែលដុះមេរោគ"
$data[32,0] = "Generate .acorn file"
$data[32,1] = "បង្កើតឯកសារ .acorn"
$data[33,0] = "Generate and load .acorn from clinical and lab data"
$data[33,1] = "បង្កើត​និងប្រតិបត្តិការឯកសារ .acorn ពីទិន្នន័យគ្លីនិកនិងមន្ទីរពិសោធន៍"
$data[34,0] = "Generating .acorn"
$data[34,1] = "កំពុងបង្កើតឯកសារ .acorn"
$data[35,0] = "Get data from REDCap"
$data[35,1] = "ទាយទិន្នន័យពី REDCap"
$data[36,0] = "Get the latest production release"
$data[36,1] = "ទាញយកផលិតកម្មដែលចេញចុងក្រោយបង្អស់"
$data[37,0] = "Growth / No Growth"
$data[37,1] = "ដុះ / មិនដុះ"
$data[38,0] = "HAI point prevalence by "
$data[38,1] = "TBT"
$data[39,0] = "HAI Prevalence"
$data[39,1] = "ប្រេវ៉ានឡង់ការបង្ករោគដែលមានប្រភពពីមន្ទីរពេទ្យ"
$data[40,0] = "Horizontal bars show the size of a set of SR results while vertical bars show the number of resistant isolates for the corresponding antibiotic."
$data[40,1] = "របារផ្តេកបង្ហាញពីទំហំសំណុំនៃលទ្ធផល SR ខណៈពេលដែលរបារបញ្ឈរបង្ហាញចំនួននៃ isolates ដែលសុំានឹងឱសថអង់ទីប៊ីយ៉ូទិកដែលឆ្លើយតប។"
$data[41,0] = "Info on loaded .acorn"
$data[41,1] = "ព៌ត័មាននៃការប្រតិបត្តិការ .acorn"
$data[42,0] = "Initial &amp; Final Surveillance Diagnosis"
$data[42,1] = "TBT"
$data[43,0] = "Initial and Final Surveillance Diagnosis"
$data[43,1] = "រោគវិនិច្ឆ័យដំបូង និងចុងក្រោយនៃការចូលរួមការអង្កេតតាមដាន"
$data[44,0] = "Isolates"
$data[44,1] = "Isolates"
$data[45,0] = "Issue detected with REDCap data. Please report to ACORN data managers. Until resolution, only existing .acorn files can be used."
$data[45,1] = "បានរកឃើញបញ្ហាជាមួយនឹងទិន្នន័យ REDCap។ សូមរាយការណ៍ទៅអ្នកគ្រប់គ្រងទិន្នន័យ។ រហូតដល់មានដំណោះស្រាយ មានតែឯកសារ​ .acorn ពីរមុនដែលអាចប្រើប្រាស់បាន។​"
$data[46,0] = "It might take a couple of minutes. This window will close on completion."
$data[46,1] = "វាប្រហែលចំណាយពេលពីរបីនាទី។  ផ្ទាំងនេះនឹងបិទពេលដំណើរការចប់។"
$data[47,0] = "Lab data not provided"
$data[47,1] = "TBT"
$data[48,0] = "Lab data successfully processed!"
$data[48,1] = "ទិន្នន័យមន្ទីរពិសោធន៍ដំណើរការដោយជោគជ័យ!"
$data[49,0] = "Lab data successfully provided."
$data[49,1] = "ទិន្នន័យមន្ទីរពិសោធន៍ត្រូវបានផ្ដល់ដោយជោគជ័យ។"
$data[50,0] = "Lab dataset contains the minimal columns."
$data[50,1] = "សំណុំទិន្នន័យមន្ទីរពិសោធន៍មានជួរឈរតិច។"
$data[51,0] = "Lab dataset does not contains the minimal columns."
$data[51,1] = "សំណុំទិន្នន័យមន្ទីរពិសោធន៍មិនមានជួរឈរតិចទេ។"
$data[52,0] = "Language"
$data[52,1] = "ភាសា"
$data[53,0] = "Load .acorn"
$data[53,1] = "ប្រតិបត្តិការ .acorn"
$data[54,0] = "Load .acorn from cloud"
$data[54,1] = "ប្រតិបត្តិការ .acorn ពី​ cloud"
$data[55,0] = "Load .acorn from local file"
$data[55,1] = "ប្រតិបត្តិការ .acorn ពីឯកសារក្នុងកុំព្យូទ័រ"
$data[56,0] = "Load selected .acorn"
$data[56,1] = "ប្រតិបត្តិការ .acorn ដែលបានជ្រើសរើស"
$data[57,0] = "Loading data."
$data[57,1] = "ទិន្នន័យកំពុងតែប្រតិបត្តិការ។"
$data[58,0] = "Log in"
$data[58,1] = "ចុះឈ្មោះចូល"
$data[59,0] = "Microbiology"
$data[59,1] = "មីក្រូជីវសាស្ត្រ"
$data[60,0] = "Most frequent 10 organisms in the plot and complete listing in the table. Contaminants are in red."
$data[60,1] = "មេរោគចំនួន១០ ដែលបានរកឃើញញឹកញាប់បំផុត គឺដាក់នៅក្នុងគំនូរតាង​ និងបញ្ជីពេញលេញនៅក្នុងតារាង។ Contaminants គឺពណ៌ ក្រហម ។"
$data[61,0] = "No .acorn data loaded."
$data[61,1] = "គ្មានទិន្នន័យ .acorn ត្រូវបានផ្ទុកទេ។​"
$data[62,0] = "No .acorn has been generated"
$data[62,1] = "TBT"
$data[63,0] = "No .acorn has been saved"
$data[63,1] = "TBT"
$data[64,0] = "No Blood Culture"
$data[64,1] = "មិនមានការបណ្តុះមេរោគក្នុងឈាម"
$data[65,0] = "Not connected to internet."
$data[65,1] = "មិនបានភ្ជាប់ទៅនឹងអ៊ីនធឺណិត"
$data[66,0] = "Number of specimens per specimen type"
$data[66,1] = "ចំនួនវត្ថុវិភាគតាមប្រភេទនៃវត្ថុវិភាគ"
$data[67,0] = "Occupancy rate per type of ward per month"
$data[67,1] = "អត្រាសម្រាកព្យាបាលទៅតាមប្រភេទនីមួយៗនៃអគារក្នុងមួយខែ"
$data[68,0] = "of blood cultures grew a potential contaminant."
$data[68,1] = "នៃការបណ្តុះមេរោគក្នុងឈាមដែលដុះ potential contaminant"
$data[69,0] = "of cultures have growth."
$data[69,1] = "នៃការបណ្តុះមេរោគដែលដុះមេរោគ"
$data[70,0] = "of enrolments with blood culture."
$data[70,1] = "នៃការចូលរួមអង្កេតតាមដាន ជាមួយការបណ្តុះមេរោគក្នុងឈាម។"
$data[71,0] = "of Target Pathogens"
$data[71,1] = "នៃមេរោគគោលដៅ"
$data[72,0] = "Only isolates that have been tested against all of the drugs are included in the upset plot."
$data[72,1] = "TBT"
$data[73,0] = "Overview"
$data[73,1] = "ទិដ្ឋភាពទូទៅ"
$data[74,0] = "Password"
$data[74,1] = "ពាក្យសម្ងាត់"
$data[75,0] = "Patient Age Distribution"
$data[75,1] = "របាយអាយុអ្នកជំងឺ"
$data[76,0] = "Patient Comorbidities"
$data[76,1] = "ស្ថានភាពជំងឺបន្សំនានារបស់អ្នកជំងឺ"
$data[77,0] = "Patient enrolments"
$data[77,1] = "អ្នកជំងឺចូលរួមការអង្កេតតាមដាន"
$data[78,0] = "Patients Transferred"
$data[78,1] = "អ្នកជំងឺដែលបានបញ្ជូនពីមន្ទីរពេទ្យផ្សេង"
$data[79,0] = "Please log in"
$data[79,1] = "សូមចុះឈ្មោះចូល"
$data[80,0] = "Processing lab data."
$data[80,1] = "កំពុងដំណើរការទិន្នន័យមន្ទីរពិសោធន៍។"
$data[81,0] = "Reading lab data."
$data[81,1] = "កំពុងពិនិត្យនិងផ្ទៀងផ្ទាត់ទិន្នន័យមន្ទីរពិសោធន៍។​"
$data[82,0] = "Remove 'Not Cultured' specimens"
$data[82,1] = "TBT"
$data[83,0] = "Remove blood culture contaminants from the following visualizations"
$data[83,1] = "លុបចោលនូវការបណ្តុះមេរោគក្នុងឈាមដែល contaminants ពីរូបភាពខាងក្រោម"
$data[84,0] = "Reset Enrolments Filters"
$data[84,1] = "កំណត់តម្រងការចូលរួមអង្កេតតាមដានឡើងវិញ"
$data[85,0] = "Resistance to 3rd gen. Cephalosporins Over Time"
$data[85,1] = "សុំាទៅនឹង  Cephalosporins ជំនាន់ទី៣ Over Time"
$data[86,0] = "Resistance to Carbapenems Over Time"
$data[86,1] = "សុំាទៅនឹង Carbapenems Over Time"
$data[87,0] = "Resistance to Fluoroquinolones Over Time"
$data[87,1] = "សុំាទៅនឹង Fluoroquinolones Over Time"
$data[88,0] = "Resistance to Oxacillin Over Time"
$data[88,1] = "សុំាទៅនឹង Oxacillin Over Time"
$data[89,0] = "Resistance to Penicillin G - meningitis Over Time"
$data[89,1] = "សុំាទៅនឹង Penicillin G - meningitis Over Time"
$data[90,0] = "Resistance to Penicillin G Over Time"
$data[90,1] = "សុំាទៅនឹង Penicillin G Over Time"
$data[91,0] = "Retriving data from REDCap server."
$data[91,1] = "ទាញយកទិន្នន័យពីម៉ាស៊ីនមេ REDCap។"
$data[92,0] = "Save .acorn file"
$data[92,1] = "រក្សាទុកឯកសារ .acorn"
$data[93,0] = "Save acorn data"
$data[93,1] = "រក្សាទុកទិន្នន័យ acorn"
$data[94,0] = "Save on Server"
$data[94,1] = "រក្សាទុកក្នុងម៉ាស៊ីនមេ"
$data[95,0] = "See Breakdown by Ward"
$data[95,1] = "មើលតាមអគារនីមួយៗ"
$data[96,0] = "See by Week"
$data[96,1] = "មើលតាមសប្តាហ៍"
$data[97,0] = "Show antibiotics combinations"
$data[97,1] = "បង្ហាញការរួមបញ្ជូលគ្នានៃឱសថអង់ទីប៊ីយ៉ូទិក"
$data[98,0] = "Show comorbidities combinations"
$data[98,1] = "បង្ហាញការរួមបញ្ទូលគ្នានៃស្ថានភាពជំងឺបន្សំនានា"
$data[99,0] = "SIR Evaluation"
$data[99,1] = "ការវាយតម្លៃ SIR"
$data[100,0] = "Some D28 records (F04) don't have a matching patient enrolment (F01)."
$data[100,1] = "កំណត់ត្រា D28 (F04) ខ្លះមិនត្រូវគ្នាជាមួយនឹងអ្នកជំងឺដែលចូលរួមការអង្កេតតាមដាន (F01)។"
$data[101,0] = "Some dates of enrolment for HAI patients do have a matching date in the HAI survey dataset"
$data[101,1] = "កាលបរិច្ឆេទមួយចំនួននៃអ្នកជំងឺដែលចូលរួមការអង្កេតតាមដាន HAI មានកាលបរិច្ឆេទត្រូវគ្នាទៅនឹងសំណុំទិន្នន័យអង្កេតតាមដាន HAI"
$data[102,0] = "Some hospital outcome records (F03) don't have a matching infection episode (F02). These records have been removed."
$data[102,1] = "កំណត់ត្រាលទ្ធផលនៅមន្ទីរពេទ្យមួយចំនួន (F03) មិនត្រូវគ្នាជាមួយវគ្គការបង្ករោគ (F02)។ កំណត់ត្រាទាំងនេះ ត្រូវបានដកចេញ ។"
$data[103,0] = "Some hospital outcome records (F03) don't have a matching patient enrolment (F01)."
$data[103,1] = "កំណត់ត្រាលទ្ធផលនៅមន្ទីរពេទ្យមួយចំនួន (F03) មិនត្រូវគ្នាជាមួយនឹងអ្នកជំងឺដែលចូលរួមការអង្កេតតាមដាន (F01)។"
$data[104,0] = "Some infection episode records (F02) don't have a matching patient enrolment (F01). These records have been removed."
$data[104,1] = "កំណត់ត្រាវគ្គការបង្ករោគមួយចំនួន (F02) មិនត្រូវគ្នាជាមួយនឹងអ្នកជំងឺដែលចូលរួមការអង្កេតតាមដាន (F01)។ កំណត់ត្រាទាំងនេះ ត្រូវបានដកចេញ ។"
$data[105,0] = "Some records with a missing ACORN ID. These records have been removed."
$data[105,1] = "កំណត់ត្រាមួយចំនួនមិនមាន ACORN ID។ កំណត់ត្រាទាំងនេះបានត្រូវលុបចោល។"
$data[106,0] = "Specimen Types"
$data[106,1] = "ប្រភេទវត្ថុវិភាគ"
$data[107,0] = "Specimens"
$data[107,1] = "TBT"
$data[108,0] = "Specimens Collected"
$data[108,1] = "វត្ថុវិភាគដែលប្រមូលបាន"
$data[109,0] = "specimens per enrolment"
$data[109,1] = "វត្ថុវិភាគនៃអ្នកចូលរួមអង្កេតតាមដានម្នាក់"
$data[110,0] = "Successfully combined clinical and lab data into .acorn file"
$data[110,1] = "ទិន្នន័យគ្លីនិកនិង​​មន្ទីរពិសោធន៍ ត្រូវបានរួមបញ្ចូលដោយជោគជ័យទៅក្នុងឯកសារ .acorn"
$data[111,0] = "Successfully loaded data."
$data[111,1] = "ទិន្នន័យប្រតិបត្តិការដោយជោគជ័យ។"
$data[112,0] = "Successfully logged in."
$data[112,1] = "ការចុះឈ្មោះចូលបានជោគជ័យ។"
$data[113,0] = "Successfully saved .acorn file in the cloud. You can now explore acorn data."
$data[113,1] = "ឯកសារ .acorn ត្រូវបានរក្សាទុកក្នុង cloud ដោយជោគជ័យ។ ឥឡូវ​នេះអ្នកអាចពិនិត្យមើលទិន្នន័យ acorn។"
$data[114,0] = "Successfully saved .acorn file locally."
$data[114,1] = "ឯកសារ .acorn ត្រូវបានរក្សាទុកក្នុងកុំព្យូទ័រដោយជោគជ័យ។"
$data[115,0] = "Supply first valid clinical and lab data."
$data[115,1] = "ផ្តល់ទិន្នន័យគ្លីនិក​និងមន្ទីពិសោធន៍ដែលមានសុពលភាពជាមុន។"
$data[116,0] = "Susceptible and Intermediate are always combined in this visualisation of co-resistances."
$data[116,1] = "Susceptible and Intermediate​ គឺតែងតែរួមបញ្ជូលគ្នានៅក្នុង​គំនូសតាងនៃសហ-ភាពសុំា"
$data[117,0] = "The 10 most common initial-final diagnosis combinations:"
$data[117,1] = "ការរួមបញ្ចូលគ្នានៃរោគវិនិច្ឆ័យដំបូងនិងចុងក្រោយទាំង១០ដែលកើតច្រើនជាងគេ៖"
$data[118,0] = "The following 'patient id' are atypical cases (one HCAI/CAI with early HAI but no overlap):"
$data[118,1] = "patient id' ខាងក្រោមគឺជាករណីមិនប្រក្រតី (HCAI/CAI មួយ ជាមួយ HAI ដំបូងប៉ុន្តែមិនមានការជាន់គ្នា)៖"
$data[119,0] = "The following 'patient id' are problem case (overlapping specimen collection windows):"
$data[119,1] = "'patient id' ខាងក្រោមគឺមានបញ្ហា​(មានការជាន់គ្នាក្នុងផ្ទាំងនៃការយកវត្ថុវិភាគ)៖"
$data[120,0] = "The REDCap dataset is empty/in wrong format. Please contact ACORN support."
$data[120,1] = "សំណុំទិន្នន័យ REDCap គឺមិនមានទិន្នន័យ / ខុសទម្រង់។ សូមទាក់ទងទៅកាន់ផ្នែកគ្រាំទ្រ ACORN។"
$data[121,0] = "The REDCap dataset is in the right format."
$data[121,1] = "សំណុំទិន្នន័យ REDCap គឺត្រឹមត្រូវតាមទម្រង់។"
$data[122,0] = "There are D28 follow-up done before the expected D28 date."
$data[122,1] = "មានការតាមដានបន្ត D28 បំពេញមុនថ្ងៃ D28 ដែលរំពឹងទុក។"
$data[123,0] = "There are multiple F02 with identical ACORN ID, admission date, and episode enrolment date."
$data[123,1] = "មាន F02 ច្រើនដែលមាន ACORN ID កាលបរិច្ឆេទចូលសម្រាកពេទ្យ និងកាលបរិច្ឆេទចូលរួមវគ្គអង្កេតតាមដានដូចគ្នា។"
$data[124,0] = "There are no atypical case (one HCAI/CAI with early HAI but no overlap)."
$data[124,1] = "មិនមានករណីមិនធម្មតាទេ (HCAI/CAI មួយ ជាមួយ HAI ដំបូងប៉ុន្តែមិនមានការជាន់គ្នា)"
$data[125,0] = "There are no D28 follow-up done before the expected D28 date."
$data[125,1] = "គ្មានការតាមដានបន្ត D28 បំពេញមុនថ្ងៃ D28 ដែលរំពឹងទុកនោះទេ។"
$data[126,0] = "There are no isolate with valid AST results. Please contact ACORN support."
$data[126,1] = "គ្មាន isolate ត្រូវជាមួយនឹងលទ្ធផល AST ទេ។ សូមទាក់ទងទៅកាន់ផ្នែកគ្រាំទ្រ ACORN។"
$data[127,0] = "There are no multiple F02 with identical ACORN ID, admission date, and episode enrolment date."
$data[127,1] = "គ្មាន F02 ដែលមាន ACORN ID កាលបរិច្ឆេទចូលសម្រាកពេទ្យ និងកាលបរិច្ឆេទចូលរួមវគ្គអង្កេតតាមដានដូចគ្នានោះទេ។"
$data[128,0] = "There are no problem case (overlapping specimen collection windows)"
$data[128,1] = "មិនមានបញ្ហា (មានការជាន់គ្នា ក្នុងផ្ទាំងនៃការយកវត្ថុវិភាគ)។"
$data[129,0] = "There are rows for which 'specdate' are after today."
$data[129,1] = "មានជួរដេកដែល 'specdate' គឺបន្ទាប់ពីថ្ងៃនេះ។"
$data[130,0] = "There are rows with missing 'orgname'."
$data[130,1] = "មានជួរដេកដែលបាត់ ‘orgname'​។"
$data[131,0] = "There are rows with missing 'patid'."
$data[131,1] = "មានជួរដេកដែលបាត់ ‘patid'។"
$data[132,0] = "There are rows with missing 'specdate'."
$data[132,1] = "មានជួរដេកដែលបាត់ ‘specdate'។"
$data[133,0] = "There are rows with missing 'specgroup'."
$data[133,1] = "មានជួរដេកដែលបាត់ ‘specgroup'។"
$data[134,0] = "There are rows with missing 'specid'."
$data[134,1] = "មានជួរដេកដែលបាត់ ‘specid'។"
$data[135,0] = "There is a critical issue with clinical data. The issue should be fixed in REDCap."
$data[135,1] = "មានបញ្ហាសំខាន់មួយក្នុងផ្នែកទិន្នន័យគ្លីនិក។ បញ្ហាគួរតែកែសម្រួលក្នុង REDCap។"
$data[136,0] = "There is no data to display for this organism."
$data[136,1] = "គ្មានទិន្នន័យ ដើម្បីបង្ហាញសម្រាប់មេរោគនេះទេ។"
$data[137,0] = "There is no HAI survey data"
$data[137,1] = "គ្មានទិន្នន័យការអង្កេតតាមដាន HAI ទេ"
$data[138,0] = "Trying to save .acorn file on server."
$data[138,1] = "កំពុងព្យាយាមរក្សាទុកឯកសារ .acorn ក្នុងម៉ាស៊ីនមេ។​"
$data[139,0] = "Updated Charlson Comorbidity Index (uCCI)"
$data[139,1] = "បច្ចុប្បន្នភាព Charlson Comorbidity Index (uCCI)"
$data[140,0] = "User"
$data[140,1] = "អ្នកប្រើប្រាស់"
$data[141,0] = "Variables in Table:"
$data[141,1] = "អថេរក្នុងតារាង៖"
$data[142,0] = "Ward Occupancy Rates"
$data[142,1] = "អត្រាសម្រាកពេទ្យតាមអគារ"
$data[143,0] = "We couldn't download the lab codes file. Please contact ACORN support."
$data[143,1] = "យើងមិនអាចទាញយក ឯកសារកូដមន្ទីរពិសោធន៍បានទេ ។​សូមទំនាក់ទំនងទៅកាន់ ACORN support។"
$data[144,0] = "We couldn't download the lab data dictionary. Please contact ACORN support"
$data[144,1] = "យើងមិនអាចទាញយកវចនានុក្រមទិន្នន័យមន្ទីរពិសោធន៍បាននោះទេ។ សូមទាក់ទងទៅកាន់ផ្នែកគាំទ្រ ACORN"
$data[145,0] = "Welcome"
$data[145,1] = "សូមស្វាគមន៍"
$data[146,0] = "What do you want to do?"
$data[146,1] = "តើអ្នកចង់ធ្វើអ្វី?"
$data[147,0] = "With Microbiology"
$data[147,1] = "ជាមួយនឹងមីក្រូជីវសាស្ត្រ"
$data[148,0] = "Wrong connection credentials."
$data[148,1] = "ព័ត៌មានសម្គាល់ការតភ្ជាប់មិនត្រឹមត្រូវ។"
$data[149,0] = "You are running ACORN dashboard"
$data[149,1] = "អ្នកកំពុងដំណើរការផ្ទាំងគ្រប់គ្រង ACORN"
$data[150,0] = "You can check here if it's the latest production release."
$data[150,1] = "អ្នកអាចពិនិត្យមើលនៅទីនេះ ថាតើផលិតកម្មជាជំនាន់ចុងក្រោយបង្អស់ឬទេ"
$data[151,0] = "Your ACORN dashboard is up to date"
$data[151,1] = "ផ្ទាំងគ្រប់គ្រង ACORN របស់អ្នកគឺជាជំនាន់ចុងក្រោយបង្អស់"
$data[152,0] = "Follow us on Twitter"
$data[152,1] = "តាមដានយើងនៅលើ Twitter"
$data[153,0] = "Records in Lab data and BSI forms:"
$data[153,1] = "កត់ត្រាក្នុងទិន្នន័យគ្លីនិកនិងមន្ទីរពិសោធន៍ក្នុងទម្រង់ BSI៖"

$ws.Range("A33:B186").Value2 = $data

Write-Output "Applied translation updates; new dimension should be A1:B186"
